# Rename the existing sheet from "Sheet1" to "Indices".
$wb = $excel.ActiveWorkbook
$wsIndices = $wb.Worksheets.Item(1)
$wsIndices.Name = "Indices"

# Add a new "FoF" worksheet right after "Indices" to hold the FoF summary table.
$wsFoF = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsIndices)
$wsFoF.Name = "FoF"

# Header.
$wsFoF.Range("B1").Value = "FoF"

# Summary statistics (label in column A, value in column B).
$wsFoF.Range("A2").Value = "Return"
$wsFoF.Range("B2").Value = 0.049341

$wsFoF.Range("A3").Value = "M squared"
$wsFoF.Range("B3").Value = 0.049341

$wsFoF.Range("A4").Value = "Volatility"
$wsFoF.Range("B4").Value = 0.056391

$wsFoF.Range("A5").Value = "MDD"
$wsFoF.Range("B5").Value = 0.222035

$wsFoF.Range("A6").Value = "CVaR"
$wsFoF.Range("B6").Value = 0.072119

$wsFoF.Range("A7").Value = "CDaR"
$wsFoF.Range("B7").Value = 0.242259

$wsFoF.Range("A8").Value = "Sharpe"
$wsFoF.Range("B8").Value = 0.122638

$wsFoF.Range("A9").Value = "Calmar"
$wsFoF.Range("B9").Value = 0.031147

$wsFoF.Range("A10").Value = "R squared"
$wsFoF.Range("B10").Value = 1

$wsFoF.Range("A11").Value = "Corr. Stocks"
$wsFoF.Range("B11").Value = 0.647835

$wsFoF.Range("A12").Value = "Corr. Bonds"
$wsFoF.Range("B12").Value = 0.360859

$wsFoF.Range("A13").Value = "Corr. FoF"
$wsFoF.Range("B13").Value = 1

$wsFoF.Range("A14").Value = "Turnover"
$wsFoF.Range("B14").Value = 0

# Match the saved view state on Indices: column K selected (scrolled right).
$wsIndices.Activate()
$wsIndices.Range("K1:K1048576").Select()

# Match the saved view state: FoF zoomed in at 175%, selection resting on C20.
# FoF is left as the active sheet/tab, same as the saved workbook.
$wsFoF.Activate()
$excel.ActiveWindow.Zoom = 175
$wsFoF.Range("C20").Select()
